# Auto-generated Excel COM-interop script
# Applies cached numeric value updates (e.g. refreshed market price data)
# to the Leve profit-calculation sheets (currentAveragePrice / LevePrice / LeveProfit columns).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7696.2085
$ws.Range("I76").Value = 10343.214
$ws.Range("J76").Value = 3990.4
$ws.Range("K76").Value = 10343.214
$ws.Range("L76").Value = 3990.4
$ws.Range("M76").Value = -10028.214
$ws.Range("N76").Value = -4620.4
$ws.Range("H79").Value = 7696.2085
$ws.Range("I79").Value = 10343.214
$ws.Range("J79").Value = 3990.4
$ws.Range("K79").Value = 10343.214
$ws.Range("L79").Value = 3990.4
$ws.Range("M79").Value = -9251.214
$ws.Range("N79").Value = -6174.4
$ws.Range("H92").Value = 569.8570999999999
$ws.Range("I92").Value = 498.16666
$ws.Range("K92").Value = 498.16666
$ws.Range("M92").Value = 749.83334
$ws.Range("H101").Value = 3514.3635
$ws.Range("J101").Value = 4583.3335
$ws.Range("L101").Value = 13750.0005
$ws.Range("N101").Value = -16994.0005
$ws.Range("H111").Value = 11115067
$ws.Range("I111").Value = 20836138
$ws.Range("J111").Value = 5271.143
$ws.Range("K111").Value = 62508414
$ws.Range("L111").Value = 15813.429
$ws.Range("M111").Value = -62505347
$ws.Range("N111").Value = -21947.429
$ws.Range("H113").Value = 35717452
$ws.Range("I113").Value = 83335416
$ws.Range("J113").Value = 3980.875
$ws.Range("K113").Value = 83335416
$ws.Range("L113").Value = 3980.875
$ws.Range("M113").Value = -83332162
$ws.Range("N113").Value = -10488.875
$ws.Range("H132").Value = 2465.8928
$ws.Range("I132").Value = 1427.8
$ws.Range("K132").Value = 4283.4
$ws.Range("M132").Value = -1753.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 38132.668
$ws.Range("J44").Value = 38132.668
$ws.Range("L44").Value = 38132.668
$ws.Range("N44").Value = -39108.668
$ws.Range("H55").Value = 25627.7
$ws.Range("J55").Value = 25627.7
$ws.Range("L55").Value = 25627.7
$ws.Range("N55").Value = -26257.7
$ws.Range("H80").Value = 25707.111
$ws.Range("J80").Value = 25707.111
$ws.Range("L80").Value = 25707.111
$ws.Range("N80").Value = -27703.111
$ws.Range("H83").Value = 25707.111
$ws.Range("J83").Value = 25707.111
$ws.Range("L83").Value = 77121.333
$ws.Range("N83").Value = -87105.333
$ws.Range("H96").Value = 28422
$ws.Range("J96").Value = 28422
$ws.Range("L96").Value = 28422
$ws.Range("N96").Value = -33914
$ws.Range("H97").Value = 2378.5334
$ws.Range("I97").Value = 1316.2693
$ws.Range("J97").Value = 3832.158
$ws.Range("K97").Value = 1316.2693
$ws.Range("L97").Value = 3832.158
$ws.Range("M97").Value = -820.2692999999999
$ws.Range("N97").Value = -4824.157999999999
$ws.Range("H132").Value = 3033.2083
$ws.Range("I132").Value = 2996.2632
$ws.Range("J132").Value = 3173.6
$ws.Range("K132").Value = 8988.7896
$ws.Range("L132").Value = 9520.799999999999
$ws.Range("M132").Value = -6458.7896
$ws.Range("N132").Value = -14580.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2966.125
$ws.Range("I3").Value = 1247.6
$ws.Range("K3").Value = 1247.6
$ws.Range("M3").Value = -1133.6
$ws.Range("H35").Value = 33125.332
$ws.Range("J35").Value = 33125.332
$ws.Range("L35").Value = 33125.332
$ws.Range("N35").Value = -33745.332
$ws.Range("H36").Value = 1150
$ws.Range("I36").Value = 1150
$ws.Range("K36").Value = 1150
$ws.Range("M36").Value = -616

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1970.7
$ws.Range("I31").Value = 1357.08
$ws.Range("J31").Value = 5038.8
$ws.Range("K31").Value = 1357.08
$ws.Range("L31").Value = 5038.8
$ws.Range("M31").Value = -1062.08
$ws.Range("N31").Value = -5628.8
$ws.Range("H34").Value = 1970.7
$ws.Range("I34").Value = 1357.08
$ws.Range("J34").Value = 5038.8
$ws.Range("K34").Value = 1357.08
$ws.Range("L34").Value = 5038.8
$ws.Range("M34").Value = -1155.08
$ws.Range("N34").Value = -5442.8
$ws.Range("H35").Value = 8980.299999999999
$ws.Range("I35").Value = 2600.4285
$ws.Range("J35").Value = 23866.666
$ws.Range("K35").Value = 2600.4285
$ws.Range("L35").Value = 23866.666
$ws.Range("M35").Value = -2306.4285
$ws.Range("N35").Value = -24454.666
$ws.Range("H58").Value = 1219.8684
$ws.Range("I58").Value = 1387.8
$ws.Range("J58").Value = 1033.2778
$ws.Range("K58").Value = 1387.8
$ws.Range("L58").Value = 1033.2778
$ws.Range("M58").Value = -1184.8
$ws.Range("N58").Value = -1439.2778
$ws.Range("H105").Value = 3445.5454
$ws.Range("I105").Value = 1651.6666
$ws.Range("K105").Value = 1651.6666
$ws.Range("M105").Value = 95.33339999999998
$ws.Range("H136").Value = 1219.8684
$ws.Range("I136").Value = 1387.8
$ws.Range("J136").Value = 1033.2778
$ws.Range("K136").Value = 4163.4
$ws.Range("L136").Value = 3099.8334
$ws.Range("M136").Value = -1613.4
$ws.Range("N136").Value = -8199.8334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 8641
$ws.Range("I82").Value = 346.5
$ws.Range("J82").Value = 10299.9
$ws.Range("K82").Value = 1039.5
$ws.Range("L82").Value = 30899.7
$ws.Range("M82").Value = -633.5
$ws.Range("N82").Value = -31711.7
$ws.Range("H85").Value = 8641
$ws.Range("I85").Value = 346.5
$ws.Range("J85").Value = 10299.9
$ws.Range("K85").Value = 1039.5
$ws.Range("L85").Value = 30899.7
$ws.Range("M85").Value = 364.5
$ws.Range("N85").Value = -33707.7
$ws.Range("H132").Value = 982
$ws.Range("I132").Value = 476
$ws.Range("K132").Value = 4284
$ws.Range("M132").Value = -1754

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 38226
$ws.Range("J123").Value = 38226
$ws.Range("L123").Value = 38226
$ws.Range("N123").Value = -43126

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2475
$ws.Range("J40").Value = 3900
$ws.Range("L40").Value = 3900
$ws.Range("N40").Value = -4172
$ws.Range("H68").Value = 2666.5557
$ws.Range("I68").Value = 1760
$ws.Range("J68").Value = 3799.75
$ws.Range("K68").Value = 1760
$ws.Range("L68").Value = 3799.75
$ws.Range("M68").Value = -1011
$ws.Range("N68").Value = -5297.75
$ws.Range("H71").Value = 2666.5557
$ws.Range("I71").Value = 1760
$ws.Range("J71").Value = 3799.75
$ws.Range("K71").Value = 8800
$ws.Range("L71").Value = 18998.75
$ws.Range("M71").Value = -5056
$ws.Range("N71").Value = -26486.75
$ws.Range("H132").Value = 3184.1
$ws.Range("I132").Value = 2309
$ws.Range("J132").Value = 4059.2
$ws.Range("K132").Value = 6927
$ws.Range("L132").Value = 12177.6
$ws.Range("M132").Value = -4397
$ws.Range("N132").Value = -17237.6
$ws.Range("H133").Value = 34397.855
$ws.Range("J133").Value = 34397.855
$ws.Range("L133").Value = 34397.855
$ws.Range("N133").Value = -39457.855

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3155.6
$ws.Range("I132").Value = 3214.0833
$ws.Range("J132").Value = 3067.875
$ws.Range("K132").Value = 9642.249899999999
$ws.Range("L132").Value = 9203.625
$ws.Range("M132").Value = -7112.249899999999
$ws.Range("N132").Value = -14263.625

